$d = $word.ActiveDocument

# The section being rewritten spans the paragraph starting with "UP -
# Subdivision" through the paragraph ending with the old "Total Amount:
# ... $4,760.00". Locate it by content rather than trusting a hardcoded
# index, then replace the whole block in one shot via InsertXML (which
# replaces the exact Range's contents).
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t.StartsWith("UP")) {
        $startPara = $p
    }
    if ($t.StartsWith("Total Amount:")) {
        $endPara = $p
        break
    }
}

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Supporting Documentation for Invoice #15________</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t>UP –___________Subdivision</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Reference #_______, MP(s) ______</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Location:  __________ County, State</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>__</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:t>__</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">/2015 &amp; </w:t>
      </w:r>
      <w:r>
        <w:t>___</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:t>__</w:t>
      </w:r>
      <w:r>
        <w:t>/2015 - Dates Tribal Fees were paid</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Tr</w:t>
      </w:r>
      <w:r>
        <w:t>ibal Consultation Fees for TCNS</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ________</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tribe</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Fee</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">   </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:t>$500.00</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>GSS Administration Fee</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>$  40</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>.00</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Total Amount:</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:r>
        <w:t>$</w:t>
      </w:r>
      <w:r>
        <w:t>54</w:t>
      </w:r>
      <w:r>
        <w:t>0.00</w:t>
      </w:r>
    </w:p>
'@

[void]$r.InsertXML($xml)
